$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.028899120726830971
$ws.Range("B1").Value = -0.028899122359886169

$ws.Range("A2").Value = 0.015873402094639312
$ws.Range("B2").Value = -0.015873403805138374

$ws.Range("A3").Value = -0.051037648898651193
$ws.Range("B3").Value = 0.051037647263369755

$ws.Range("A4").Value = -0.054427067325529835
$ws.Range("B4").Value = 0.054427065649962976

$ws.Range("A5").Value = 0.025159800086637662
$ws.Range("B5").Value = -0.025159801829061825
